# Auto-generated update of betting odds cells per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 15).Value = 1.17
$ws.Cells.Item(2, 16).Value = 5

# Row 9
$ws.Cells.Item(9, 7).Value = 2.15
$ws.Cells.Item(9, 9).Value = 3.25
$ws.Cells.Item(9, 10).Value = 2.88
$ws.Cells.Item(9, 12).Value = 3.75
$ws.Cells.Item(9, 15).Value = 1.25
$ws.Cells.Item(9, 16).Value = 3.75
$ws.Cells.Item(9, 17).Value = 1.93
$ws.Cells.Item(9, 18).Value = 1.93
$ws.Cells.Item(9, 24).Value = 11
$ws.Cells.Item(9, 26).Value = 21
$ws.Cells.Item(9, 30).Value = 6.5
$ws.Cells.Item(9, 50).Value = 5
$ws.Cells.Item(9, 51).Value = 17
$ws.Cells.Item(9, 55).Value = 151

# Row 10
$ws.Cells.Item(10, 7).Value = 2.75
$ws.Cells.Item(10, 8).Value = 3.25
$ws.Cells.Item(10, 9).Value = 2.27
$ws.Cells.Item(10, 13).Value = 1.02
$ws.Cells.Item(10, 14).Value = 10.8
$ws.Cells.Item(10, 15).Value = 1.26
$ws.Cells.Item(10, 16).Value = 3.22
$ws.Cells.Item(10, 17).Value = 1.9
$ws.Cells.Item(10, 18).Value = 1.72
$ws.Cells.Item(10, 19).Value = 1.36
$ws.Cells.Item(10, 20).Value = 2.52
$ws.Cells.Item(10, 21).Value = 1.72
$ws.Cells.Item(10, 22).Value = 2.07
$ws.Cells.Item(10, 23).Value = 7.3
$ws.Cells.Item(10, 24).Value = 11.25
$ws.Cells.Item(10, 25).Value = 8.75
$ws.Cells.Item(10, 27).Value = 19
$ws.Cells.Item(10, 28).Value = 26
$ws.Cells.Item(10, 29).Value = 9.25
$ws.Cells.Item(10, 30).Value = 5.6
$ws.Cells.Item(10, 31).Value = 11.75
$ws.Cells.Item(10, 32).Value = 50
$ws.Cells.Item(10, 33).Value = 350
$ws.Cells.Item(10, 34).Value = 6.7
$ws.Cells.Item(10, 35).Value = 9.25
$ws.Cells.Item(10, 36).Value = 7.8
$ws.Cells.Item(10, 37).Value = 18
$ws.Cells.Item(10, 38).Value = 15
$ws.Cells.Item(10, 39).Value = 23
$ws.Cells.Item(10, 40).Value = 4.7
$ws.Cells.Item(10, 41).Value = 15
$ws.Cells.Item(10, 42).Value = 23
$ws.Cells.Item(10, 43).Value = 70
$ws.Cells.Item(10, 44).Value = 100
$ws.Cells.Item(10, 45).Value = 300
$ws.Cells.Item(10, 46).Value = 2.52
$ws.Cells.Item(10, 47).Value = 7.1
$ws.Cells.Item(10, 48).Value = 65
$ws.Cells.Item(10, 50).Value = 4.2
$ws.Cells.Item(10, 51).Value = 11.75
$ws.Cells.Item(10, 52).Value = 20

# Row 11
$ws.Cells.Item(11, 15).Value = 1.26
$ws.Cells.Item(11, 16).Value = 3.22

# Row 13
$ws.Cells.Item(13, 7).Value = 2.67
$ws.Cells.Item(13, 8).Value = 3.15
$ws.Cells.Item(13, 9).Value = 2.5
$ws.Cells.Item(13, 10).Value = 3.2
$ws.Cells.Item(13, 11).Value = 2.12
$ws.Cells.Item(13, 12).Value = 3
$ws.Cells.Item(13, 14).Value = 8.2
$ws.Cells.Item(13, 23).Value = 8.25
$ws.Cells.Item(13, 24).Value = 13.5
$ws.Cells.Item(13, 25).Value = 10
$ws.Cells.Item(13, 26).Value = 30
$ws.Cells.Item(13, 27).Value = 23
$ws.Cells.Item(13, 30).Value = 6.1
$ws.Cells.Item(13, 34).Value = 8.25
$ws.Cells.Item(13, 35).Value = 12.5
$ws.Cells.Item(13, 36).Value = 9.5
$ws.Cells.Item(13, 37).Value = 27
$ws.Cells.Item(13, 38).Value = 21
$ws.Cells.Item(13, 40).Value = 4.6
$ws.Cells.Item(13, 41).Value = 14
$ws.Cells.Item(13, 43).Value = 60
$ws.Cells.Item(13, 46).Value = 2.8
$ws.Cells.Item(13, 47).Value = 6.6
$ws.Cells.Item(13, 50).Value = 4.45
$ws.Cells.Item(13, 51).Value = 13
$ws.Cells.Item(13, 52).Value = 19
$ws.Cells.Item(13, 53).Value = 50
$ws.Cells.Item(13, 55).Value = 200

# Row 14
$ws.Cells.Item(14, 7).Value = 3.3
$ws.Cells.Item(14, 8).Value = 3.1
$ws.Cells.Item(14, 9).Value = 2.15
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 2.75
$ws.Cells.Item(14, 14).Value = 7.6
$ws.Cells.Item(14, 20).Value = 2.45
$ws.Cells.Item(14, 22).Value = 2
$ws.Cells.Item(14, 24).Value = 17.5
$ws.Cells.Item(14, 26).Value = 45
$ws.Cells.Item(14, 27).Value = 30
$ws.Cells.Item(14, 28).Value = 35
$ws.Cells.Item(14, 29).Value = 9.25
$ws.Cells.Item(14, 30).Value = 6.1
$ws.Cells.Item(14, 31).Value = 12.5
$ws.Cells.Item(14, 34).Value = 8.25
$ws.Cells.Item(14, 35).Value = 11.25
$ws.Cells.Item(14, 37).Value = 22
$ws.Cells.Item(14, 39).Value = 24
$ws.Cells.Item(14, 46).Value = 2.42
$ws.Cells.Item(14, 47).Value = 6.8
$ws.Cells.Item(14, 51).Value = 11.25

# Row 16
$ws.Cells.Item(16, 17).Value = 1.7
$ws.Cells.Item(16, 18).Value = 2.1

# Row 17
$ws.Cells.Item(17, 7).Value = 1.44
$ws.Cells.Item(17, 8).Value = 4.33
$ws.Cells.Item(17, 9).Value = 7
$ws.Cells.Item(17, 10).Value = 1.95
$ws.Cells.Item(17, 17).Value = 1.75
$ws.Cells.Item(17, 18).Value = 2.05
$ws.Cells.Item(17, 19).Value = 1.33
$ws.Cells.Item(17, 20).Value = 3.25
$ws.Cells.Item(17, 30).Value = 8.5
$ws.Cells.Item(17, 37).Value = 81
$ws.Cells.Item(17, 41).Value = 7
$ws.Cells.Item(17, 46).Value = 3.25

# Row 18
$ws.Cells.Item(18, 7).Value = 1.17
$ws.Cells.Item(18, 10).Value = 1.53
$ws.Cells.Item(18, 13).Value = 1.03
$ws.Cells.Item(18, 14).Value = 17
$ws.Cells.Item(18, 23).Value = 8.5
$ws.Cells.Item(18, 34).Value = 41
$ws.Cells.Item(18, 39).Value = 101
$ws.Cells.Item(18, 41).Value = 5
$ws.Cells.Item(18, 54).Value = 301

# Row 19
$ws.Cells.Item(19, 7).Value = 3.7
$ws.Cells.Item(19, 9).Value = 2.05
$ws.Cells.Item(19, 10).Value = 4.33
$ws.Cells.Item(19, 12).Value = 2.75
$ws.Cells.Item(19, 14).Value = 8.5
$ws.Cells.Item(19, 21).Value = 1.95
$ws.Cells.Item(19, 22).Value = 1.8
$ws.Cells.Item(19, 24).Value = 19
$ws.Cells.Item(19, 27).Value = 34
$ws.Cells.Item(19, 30).Value = 6.5
$ws.Cells.Item(19, 33).Value = 351
$ws.Cells.Item(19, 34).Value = 6.5
$ws.Cells.Item(19, 35).Value = 9
$ws.Cells.Item(19, 38).Value = 17

# Row 20
$ws.Cells.Item(20, 33).Value = 1250

# Row 22
$ws.Cells.Item(22, 7).Value = 1.47
$ws.Cells.Item(22, 8).Value = 4.2
$ws.Cells.Item(22, 9).Value = 5.4
$ws.Cells.Item(22, 10).Value = 1.98
$ws.Cells.Item(22, 11).Value = 2.4
$ws.Cells.Item(22, 12).Value = 5.3
$ws.Cells.Item(22, 13).Value = 1.04
$ws.Cells.Item(22, 14).Value = 8.75
$ws.Cells.Item(22, 15).Value = 1.19
$ws.Cells.Item(22, 16).Value = 4.1
$ws.Cells.Item(22, 17).Value = 1.6
$ws.Cells.Item(22, 18).Value = 2.2
$ws.Cells.Item(22, 19).Value = 1.3
$ws.Cells.Item(22, 20).Value = 3.2
$ws.Cells.Item(22, 21).Value = 1.72
$ws.Cells.Item(22, 22).Value = 2
$ws.Cells.Item(22, 24).Value = 7.7
$ws.Cells.Item(22, 26).Value = 10.75
$ws.Cells.Item(22, 27).Value = 11.25
$ws.Cells.Item(22, 29).Value = 8.75
$ws.Cells.Item(22, 30).Value = 8.5
$ws.Cells.Item(22, 31).Value = 16
$ws.Cells.Item(22, 34).Value = 18
$ws.Cells.Item(22, 35).Value = 37
$ws.Cells.Item(22, 36).Value = 17.5
$ws.Cells.Item(22, 37).Value = 110
$ws.Cells.Item(22, 38).Value = 50
$ws.Cells.Item(22, 39).Value = 45
$ws.Cells.Item(22, 40).Value = 3.5
$ws.Cells.Item(22, 41).Value = 6.9
$ws.Cells.Item(22, 42).Value = 15
$ws.Cells.Item(22, 43).Value = 19.5
$ws.Cells.Item(22, 44).Value = 45
$ws.Cells.Item(22, 46).Value = 3.2
$ws.Cells.Item(22, 50).Value = 7.3
$ws.Cells.Item(22, 51).Value = 30
$ws.Cells.Item(22, 52).Value = 30
$ws.Cells.Item(22, 53).Value = 175
$ws.Cells.Item(22, 54).Value = 175

# Row 23
$ws.Cells.Item(23, 7).Value = 6.5
$ws.Cells.Item(23, 9).Value = 1.42
$ws.Cells.Item(23, 10).Value = 6
$ws.Cells.Item(23, 12).Value = 1.91
$ws.Cells.Item(23, 15).Value = 1.17
$ws.Cells.Item(23, 16).Value = 5
$ws.Cells.Item(23, 17).Value = 1.53
$ws.Cells.Item(23, 18).Value = 2.4
$ws.Cells.Item(23, 21).Value = 1.67
$ws.Cells.Item(23, 22).Value = 2.1
$ws.Cells.Item(23, 23).Value = 21
$ws.Cells.Item(23, 25).Value = 19
$ws.Cells.Item(23, 26).Value = 67
$ws.Cells.Item(23, 37).Value = 11
$ws.Cells.Item(23, 40).Value = 8
$ws.Cells.Item(23, 41).Value = 29
$ws.Cells.Item(23, 43).Value = 101
$ws.Cells.Item(23, 44).Value = 101
$ws.Cells.Item(23, 45).Value = 151
$ws.Cells.Item(23, 47).Value = 8
$ws.Cells.Item(23, 51).Value = 7
$ws.Cells.Item(23, 53).Value = 19
$ws.Cells.Item(23, 54).Value = 41

# Row 25
$ws.Cells.Item(25, 7).Value = 1.29
$ws.Cells.Item(25, 8).Value = 4.85
$ws.Cells.Item(25, 10).Value = 1.75
$ws.Cells.Item(25, 11).Value = 2.37
$ws.Cells.Item(25, 12).Value = 7.8
$ws.Cells.Item(25, 14).Value = 14
$ws.Cells.Item(25, 16).Value = 3.5
$ws.Cells.Item(25, 17).Value = 1.65
$ws.Cells.Item(25, 18).Value = 1.98
$ws.Cells.Item(25, 23).Value = 6.6
$ws.Cells.Item(25, 25).Value = 8.75
$ws.Cells.Item(25, 26).Value = 7.6
$ws.Cells.Item(25, 27).Value = 11.5
$ws.Cells.Item(25, 29).Value = 12
$ws.Cells.Item(25, 30).Value = 10
$ws.Cells.Item(25, 34).Value = 22
$ws.Cells.Item(25, 41).Value = 5.6
$ws.Cells.Item(25, 43).Value = 15
$ws.Cells.Item(25, 46).Value = 2.92
